$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1: text "123456789" entered with a leading apostrophe (quote-prefixed,
# stored as text even though it looks like a number)
$ws.Range("C1").Value = "'123456789"

# D1: genuine number 1234890
$ws.Range("D1").Value = 1234890

# Carry the quote-prefix cell format from C1 onto D1 (matches the source
# workbook where both cells share the same style index) without touching
# D1's numeric value. -4122 == xlPasteFormats.
$xlPasteFormats = -4122
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial($xlPasteFormats)

# Column C was resized after the new data was entered.
$ws.Columns("C:C").ColumnWidth = 9.6

# Final selection left on H6 before save.
$ws.Range("H6").Select() | Out-Null
